# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (columns H-N) for a
# handful of leve rows across the per-job Asura_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1247.85
$ws.Range("J17").Value = 1247.85
$ws.Range("L17").Value = 3743.55
$ws.Range("N17").Value = -4079.55
$ws.Range("H64").Value = 3169.9666
$ws.Range("I64").Value = 2949.9167
$ws.Range("J64").Value = 3316.6667
$ws.Range("K64").Value = 2949.9167
$ws.Range("L64").Value = 3316.6667
$ws.Range("M64").Value = -2701.9167
$ws.Range("N64").Value = -3812.6667
$ws.Range("H67").Value = 3169.9666
$ws.Range("I67").Value = 2949.9167
$ws.Range("J67").Value = 3316.6667
$ws.Range("K67").Value = 2949.9167
$ws.Range("L67").Value = 3316.6667
$ws.Range("M67").Value = -2091.9167
$ws.Range("N67").Value = -5032.6667
$ws.Range("H103").Value = 584.2222
$ws.Range("J103").Value = 596.5
$ws.Range("L103").Value = 1789.5
$ws.Range("N103").Value = -2961.5
$ws.Range("H139").Value = 61038.668
$ws.Range("J139").Value = 61038.668
$ws.Range("L139").Value = 61038.668
$ws.Range("N139").Value = -71318.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1072.75
$ws.Range("I45").Value = 956
$ws.Range("J45").Value = 1189.5
$ws.Range("K45").Value = 956
$ws.Range("L45").Value = 1189.5
$ws.Range("M45").Value = -579
$ws.Range("N45").Value = -1943.5
$ws.Range("H107").Value = 23250
$ws.Range("J107").Value = 23250
$ws.Range("L107").Value = 23250
$ws.Range("N107").Value = -30930
$ws.Range("H122").Value = 1281.6786
$ws.Range("I122").Value = 1156.4117
$ws.Range("J122").Value = 1475.2727
$ws.Range("K122").Value = 3469.2351
$ws.Range("L122").Value = 4425.8181
$ws.Range("M122").Value = -1019.2351
$ws.Range("N122").Value = -9325.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 29649.541
$ws.Range("J20").Value = 2583.5715
$ws.Range("L20").Value = 2583.5715
$ws.Range("N20").Value = -3077.5715
$ws.Range("H22").Value = 10483.667
$ws.Range("I22").Value = 15350.5
$ws.Range("K22").Value = 15350.5
$ws.Range("M22").Value = -15177.5
$ws.Range("H105").Value = 3251.2727
$ws.Range("I105").Value = 3163.9285
$ws.Range("J105").Value = 3315.6316
$ws.Range("K105").Value = 3163.9285
$ws.Range("L105").Value = 3315.6316
$ws.Range("M105").Value = -1416.9285
$ws.Range("N105").Value = -6809.6316

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1765.4921
$ws.Range("I31").Value = 1265.4667
$ws.Range("J31").Value = 3015.5557
$ws.Range("K31").Value = 1265.4667
$ws.Range("L31").Value = 3015.5557
$ws.Range("M31").Value = -970.4667
$ws.Range("N31").Value = -3605.5557
$ws.Range("H34").Value = 1765.4921
$ws.Range("I34").Value = 1265.4667
$ws.Range("J34").Value = 3015.5557
$ws.Range("K34").Value = 1265.4667
$ws.Range("L34").Value = 3015.5557
$ws.Range("M34").Value = -1063.4667
$ws.Range("N34").Value = -3419.5557
$ws.Range("H87").Value = 26680
$ws.Range("J87").Value = 26680
$ws.Range("L87").Value = 26680
$ws.Range("N87").Value = -29052
$ws.Range("H90").Value = 26680
$ws.Range("J90").Value = 26680
$ws.Range("L90").Value = 80040
$ws.Range("N90").Value = -91896
$ws.Range("H99").Value = 3362.842
$ws.Range("I99").Value = 3630.8462
$ws.Range("J99").Value = 2782.1667
$ws.Range("K99").Value = 3630.8462
$ws.Range("L99").Value = 2782.1667
$ws.Range("M99").Value = -2132.8462
$ws.Range("N99").Value = -5778.1667
$ws.Range("H126").Value = 3362.842
$ws.Range("I126").Value = 3630.8462
$ws.Range("J126").Value = 2782.1667
$ws.Range("K126").Value = 10892.5386
$ws.Range("L126").Value = 8346.5001
$ws.Range("M126").Value = -8422.5386
$ws.Range("N126").Value = -13286.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6762.222
$ws.Range("I3").Value = 5143.3335
$ws.Range("K3").Value = 15430.0005
$ws.Range("M3").Value = -15318.0005
$ws.Range("H19").Value = 5000
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 15000
$ws.Range("N19").Value = -15348
$ws.Range("H121").Value = 1178.3684
$ws.Range("J121").Value = 1277.0714
$ws.Range("L121").Value = 3831.2142
$ws.Range("N121").Value = -6451.2142
$ws.Range("H134").Value = 4548.1665
$ws.Range("I134").Value = 1306.6666
$ws.Range("J134").Value = 6349
$ws.Range("K134").Value = 3919.9998
$ws.Range("L134").Value = 19047
$ws.Range("M134").Value = 1150.0002
$ws.Range("N134").Value = -29187
$ws.Range("H136").Value = 3451.4
$ws.Range("I136").Value = 1060.4706
$ws.Range("J136").Value = 6578
$ws.Range("K136").Value = 3181.4118
$ws.Range("L136").Value = 19734
$ws.Range("M136").Value = 1918.5882
$ws.Range("N136").Value = -29934

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 63333.332
$ws.Range("J88").Value = 63333.332
$ws.Range("L88").Value = 63333.332
$ws.Range("N88").Value = -64235.332
$ws.Range("H91").Value = 63333.332
$ws.Range("J91").Value = 63333.332
$ws.Range("L91").Value = 63333.332
$ws.Range("N91").Value = -66453.332
$ws.Range("H102").Value = 5142.857
$ws.Range("I102").Value = 5000
$ws.Range("J102").Value = 5250
$ws.Range("K102").Value = 5000
$ws.Range("L102").Value = 5250
$ws.Range("M102").Value = -3378
$ws.Range("N102").Value = -8494
$ws.Range("H134").Value = 500326
$ws.Range("J134").Value = 500326
$ws.Range("L134").Value = 1500978
$ws.Range("N134").Value = -1506048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H82").Value = 1917.1428
$ws.Range("I82").Value = 1813.3572
$ws.Range("J82").Value = 2124.7144
$ws.Range("K82").Value = 1813.3572
$ws.Range("L82").Value = 2124.7144
$ws.Range("M82").Value = -1452.3572
$ws.Range("N82").Value = -2846.7144
$ws.Range("H85").Value = 1917.1428
$ws.Range("I85").Value = 1813.3572
$ws.Range("J85").Value = 2124.7144
$ws.Range("K85").Value = 1813.3572
$ws.Range("L85").Value = 2124.7144
$ws.Range("M85").Value = -565.3572
$ws.Range("N85").Value = -4620.7144
$ws.Range("H132").Value = 6999.8823
$ws.Range("I132").Value = 7312.5
$ws.Range("K132").Value = 21937.5
$ws.Range("M132").Value = -19407.5
$ws.Range("H140").Value = 27714.5
$ws.Range("J140").Value = 27714.5
$ws.Range("L140").Value = 27714.5
$ws.Range("N140").Value = -38074.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 69564.82
$ws.Range("J46").Value = 69564.82
$ws.Range("L46").Value = 69564.82
$ws.Range("N46").Value = -70026.82
$ws.Range("H93").Value = 29000
$ws.Range("J93").Value = 29000
$ws.Range("L93").Value = 29000
$ws.Range("N93").Value = -33992
$ws.Range("H123").Value = 23426.059
$ws.Range("J123").Value = 23426.059
$ws.Range("L123").Value = 23426.059
$ws.Range("N123").Value = -33226.059
$ws.Range("H134").Value = 69564.82
$ws.Range("J134").Value = 69564.82
$ws.Range("L134").Value = 208694.46
$ws.Range("N134").Value = -213764.46
